# Updating the publishing image
# Change the line dash style of four ovals on the last slide from
# "System Dash Dot" (sysDot / msoLineSquareDot) to "Dash" (dash / msoLineDash).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)

$targetNames = @("Oval 4", "Oval 23", "Oval 27", "Oval 29")

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($targetNames -contains $shp.Name) {
        $shp.Line.DashStyle = 4  # msoLineDash
    }
}
